$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly row pushes the two most recent historical rows down:
#  - old row 15 (2021-12-21 / 44551) becomes row 16, with a new, more
#    recent date (2021-12-23 / 44553) taking its place in row 15
#  - old row 16 (2021-11-19 / 44519) becomes row 17, unchanged

# Capture old row 16 values (M:T) before they are overwritten, since row 17
# needs to end up with exactly this data.
$M16old = $ws.Range("M16").Value2
$N16old = $ws.Range("N16").Value2
$O16old = $ws.Range("O16").Value2
$P16old = $ws.Range("P16").Value2
$Q16old = $ws.Range("Q16").Value2
$R16old = $ws.Range("R16").Value2
$S16old = $ws.Range("S16").Value2
$T16old = $ws.Range("T16").Value2

# Row 15: keep everything, only advance the date.
$ws.Range("D15").Value = 44553

# Row 16: keep A:C/E:L as-is, but set the date to the previous row-15 date
# and copy across the row-15 M:T figures (the old row 15 "moves" into row 16).
$ws.Range("D16").Value = 44551
$ws.Range("M16").Value = 400
$ws.Range("N16").Value = 5000
$ws.Range("O16").Value = 5500
$ws.Range("P16").Value = 5250
$ws.Range("Q16").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R16").Value = "Región del Maule"
$ws.Range("S16").Value = 3500
$ws.Range("T16").Value = 1.5

# Row 17: brand-new row holding what used to be row 16.
$ws.Range("A17").Value = 4
$ws.Range("B17").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C17").Value = "Los Lagos"
$ws.Range("D17").Value = 44519
$ws.Range("D17").NumberFormat = $ws.Range("D16").NumberFormat
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100101
$ws.Range("H17").Value = "Berries"
$ws.Range("I17").Value = 100101001
$ws.Range("J17").Value = "Arándano (blue)"
$ws.Range("K17").Value = "Sin especificar"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = $M16old
$ws.Range("N17").Value = $N16old
$ws.Range("O17").Value = $O16old
$ws.Range("P17").Value = $P16old
$ws.Range("Q17").Value = $Q16old
$ws.Range("R17").Value = $R16old
$ws.Range("S17").Value = $S16old
$ws.Range("T17").Value = $T16old
